$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (header "Förändrad") holds a date value that was bumped by one day
# (45181 -> 45182) for every data row (rows 2 through 387).
$ws.Range("C2:C387").Value = 45182
